$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.700.12'
$ws.Range("E2").Value = '  +2.21%  '

$ws.Range("D3").Value = '3.947.92'
$ws.Range("E3").Value = '  +0.89%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").NumberFormat = "General"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '519.29'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +6.46%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.16'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.61%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.622'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.16%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.997'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.08%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.730'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.22%  '

$ws.Range("E10").Value = '  +4.56%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000344'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.58%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.24'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.68%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.42'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.90%  '

$ws.Range("D14").Value = '4.570.72'
$ws.Range("E14").Value = '  +0.50%  '

$ws.Range("D15").Value = '3.954.86'
$ws.Range("E15").Value = '  +1.29%  '

$ws.Range("E16").Value = '  -0.95%  '

$ws.Range("E17").Value = '  -0.53%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.22'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +7.66%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '19.86'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.02%  '

$ws.Range("D20").Value = '69.670.31'
$ws.Range("E20").Value = '  +1.96%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '435.09'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.40%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.41'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.51%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.57'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.53%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '88.38'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.39%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.82'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.96%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.93'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +8.02%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.23'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.41%  '

$ws.Range("E28").Value = '  -4.05%  '

$ws.Range("E29").Value = '  -1.49%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '703.11'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.20%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '13.31'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.56%  '

$ws.Range("E32").Value = '  -2.67%  '

$ws.Range("E33").Value = '  -1.43%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '68.45'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +12.06%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.441'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +7.75%  '

$ws.Range("D36").Value = '0.0₃0881'
$ws.Range("E36").Value = '  +1.98%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.94'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.37%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '40.25'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.44%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.997'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.15%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.08%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0485'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.57%  '

$ws.Range("B43").Value = 'WEMIXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.11'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +6.62%  '

$ws.Range("B44").Value = 'Fetch.AI'
$ws.Range("C44").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.82'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.00%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.02'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.42%  '

$ws.Range("E46").Value = '  +1.00%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.34'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.94%  '

$ws.Range("D48").Value = '0.0₆0355'
$ws.Range("E48").Value = '  +2.64%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.00'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +6.06%  '

$ws.Range("E50").Value = '  -1.81%  '

$ws.Range("E51").Value = '  -1.97%  '
